$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.684.55"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.849.30"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'322.12"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "'1.030"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4386"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.3789"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "'0.07383"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'0.8818"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "'21.55"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.862.15"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "'5.497"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'6.699"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'0.07156"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "'84.96"
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "'1.036"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'0.000009051"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'1.030"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'15.44"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "27.692.85"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'5.288"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'11.30"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "2.085.11"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "'2.066"
$ws.Range("E25").Value = "  +7.07%  "
$ws.Range("D26").Value = "'158.82"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'18.69"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'1.988"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").Value = "'5.315"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "'117.59"
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("D31").Value = "'0.09047"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "'0.7710"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").Value = "'1.206"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'3.002"
$ws.Range("E34").Value = "  +4.31%  "
$ws.Range("D35").Value = "'4.553"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "'1.032"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "'1.150"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "'0.01970"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "'0.05266"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "'2.839"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("D41").Value = "'0.5171"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.1669"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'6.863"
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("D44").Value = "'8.700"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").Value = "'110.17"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").Value = "'10.66"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").Value = "'0.06590"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("D48").Value = "'1.032"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'1.698"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "'0.4688"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").Value = "'1.887"
$ws.Range("E51").Value = "  -0.18%  "
